$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Determine the last used row (data occupies A1:G200 -> header + 199 data rows)
$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -lt 1) { $lastRow = 1 }

for ($r = 1; $r -le $lastRow; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $cellB = $ws.Cells.Item($r, 2)

    $valA = $cellA.Value2
    $valB = $cellB.Value2

    $cellA.Value2 = $valB
    $cellB.Value2 = $valA
}
